# Day 109 and 110 - mark additional "Binary Trees" rows as done ("yes")
# with the same look-and-feel (random-ish colored fill) already used
# elsewhere on the sheet for completed items.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of target row -> a donor cell elsewhere on the sheet that already
# carries the desired cell style (colored fill + centered alignment).
# Using Copy/PasteSpecial(xlPasteFormats) reuses the existing style record
# instead of fabricating a brand new style/fill entry.
$rowStyleDonor = @{
    177 = "C22"   # style s=15
    178 = "C22"   # style s=15
    179 = "C6"    # style s=11
    180 = "C22"   # style s=15
    181 = "C22"   # style s=15
    182 = "C13"   # style s=12
    183 = "C13"   # style s=12
    184 = "C20"   # style s=14
    185 = "C13"   # style s=12
    186 = "C13"   # style s=12
    187 = "C20"   # style s=14
    188 = "C20"   # style s=14
    189 = "C13"   # style s=12
    190 = "C22"   # style s=15
    191 = "C20"   # style s=14
    192 = "C13"   # style s=12
}

foreach ($row in 177..192) {
    $donor = $rowStyleDonor[$row]
    $target = "C$row"

    $ws.Range($donor).Copy()
    $ws.Range($target).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $ws.Range($target).Value = "yes"
}

$excel.CutCopyMode = 0

# Reflect the author's final cursor position / selection on the sheet
# (Day 110's last touched cell).
$ws.Range("C178").Select()
